# Weekly fruit/vegetable price update: insert two new daily records
# (rows 216-217) for "Feria Lagunitas de Puerto Montt - Ajo", pushing
# the existing data rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 216; everything below shifts down.
$ws.Rows("216:217").Insert()

# --- New row 216 ---
$ws.Range("A216").Value2 = 4
$ws.Range("B216").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C216").Value2 = "Los Lagos"
$ws.Range("D216").Value2 = 44813
$ws.Range("E216").Value2 = 10
$ws.Range("F216").Value2 = 100112003
$ws.Range("G216").Value2 = "Ajo"
$ws.Range("H216").Value2 = "Chino"
$ws.Range("I216").Value2 = "Primera"
$ws.Range("J216").Value2 = 120
$ws.Range("K216").Value2 = 28000
$ws.Range("L216").Value2 = 28000
$ws.Range("M216").Value2 = 28000
$ws.Range("N216").Value2 = "`$/caja 10 kilos"
$ws.Range("O216").Value2 = "China"
$ws.Range("P216").Value2 = 2800
$ws.Range("Q216").Value2 = 10
$ws.Range("R216").Value2 = "Hortaliza"

# --- New row 217 ---
$ws.Range("A217").Value2 = 4
$ws.Range("B217").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C217").Value2 = "Los Lagos"
$ws.Range("D217").Value2 = 44813
$ws.Range("E217").Value2 = 10
$ws.Range("F217").Value2 = 100112003
$ws.Range("G217").Value2 = "Ajo"
$ws.Range("H217").Value2 = "Chino"
$ws.Range("I217").Value2 = "Segunda"
$ws.Range("J217").Value2 = 120
$ws.Range("K217").Value2 = 25000
$ws.Range("L217").Value2 = 25000
$ws.Range("M217").Value2 = 25000
$ws.Range("N217").Value2 = "`$/caja 10 kilos"
$ws.Range("O217").Value2 = "China"
$ws.Range("P217").Value2 = 2500
$ws.Range("Q217").Value2 = 10
$ws.Range("R217").Value2 = "Hortaliza"
